$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BD")

# New column header and Predio values, in original entry order
$ws.Range("E1").Value = "Predio"
$ws.Range("E5").Value = "CAU044_ElPrado"
$ws.Range("E18").Value = "CAU044_ElPrado_2"
$ws.Range("E17").Value = "CAU044_ElPrado_3"
$ws.Range("E19").Value = "CAU062_Beer-Sebas_1"
$ws.Range("E20").Value = "CAU062_Beer-Sebas_2"
$ws.Range("E21").Value = "CAU062_Beer-Sebas_3"
$ws.Range("E15").Value = "CAU059_Predio-Estrecho_3"
$ws.Range("E16").Value = "NAR031_Chachagui1"
$ws.Range("E23").Value = "VAL002 - Andalucia"

# Column B data updates (Etapa del proyecto)
$ws.Range("B12").Value = "Construcción"
$ws.Range("B22").Value = "Estructuración"

# Clear custom row formatting on row 21 (revert A21/B21 to default style),
# then restore C21/D21's original formatting (date format / alignment).
$ws.Rows.Item(21).ClearFormats()
$ws.Range("C2:D2").Copy()
$ws.Range("C21:D21").PasteSpecial(-4122)
$excel.CutCopyMode = $False

# Size column E to fit its new contents (matches bestFit sizing of the other columns)
$ws.Columns.Item(5).ColumnWidth = 18.63

# Apply AutoFilter over A1:E23
$ws.Range("A1:E23").AutoFilter()

# Selection/view changes
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 1
